$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = "57.991.72"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "2.470.45"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.10"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  -3.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.67"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -4.07%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.558"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  -1.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0996"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.40"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").Value = "2.910.62"
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("D14").Value = "57.926.37"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.19"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  -4.51%  "
$ws.Range("E16").Value = "  -2.37%  "
$ws.Range("D17").Value = "2.477.25"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.88"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.18"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  -2.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "320.46"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.77"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  -3.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.46"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.410"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -3.35%  "
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("E26").Value = "  -3.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.44"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  -3.15%  "
$ws.Range("D28").Value = "0.0₃0749"
$ws.Range("E28").Value = "  -4.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.40"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  -5.14%  "
$ws.Range("E30").Value = "  -4.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.07"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("E32").Value = "  -4.94%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.17"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  -1.83%  "
$ws.Range("E36").Value = "  -9.32%  "
$ws.Range("E37").Value = "  -4.38%  "
$ws.Range("E38").Value = "  -4.75%  "
$ws.Range("E39").Value = "  -3.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.49"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -4.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "276.15"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  -4.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.06"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.590"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -3.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.87"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  -3.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0912"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0493"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  -3.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0215"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  -3.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "17.11"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("D49").Value = "1.733.46"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("E50").Value = "  -1.42%  "
$ws.Range("E51").Value = "  -2.74%  "
